# Update slides for SS22.
$p = $ppt.ActivePresentation

# --- Change 1: title slide "(Group 7)" -> "(Group 8)" ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $txt = $tr.Text
            if ($txt.Contains("(Group 7)")) {
                $needle = "(Group 7)"
                $startPos = $txt.IndexOf($needle) + 1
                $found = $tr.Characters($startPos, $needle.Length)
                $found.Text = "(Group 8)"
            }
        }
    }
}

# --- Change 2: homework slide due-date run, re-dated & split into 3 runs ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $txt = $tr.Text
            if ($txt.Contains("Due date: Friday April 8 @17h")) {
                $prefix = "Due date: Friday "
                $oldDate = "April 8 "
                $newDate = "October 7 "
                $runStart = $txt.IndexOf($prefix + $oldDate + "@17h") + 1
                $oldDateStart = $runStart + $prefix.Length
                $oldDateRange = $tr.Characters($oldDateStart, $oldDate.Length)
                $oldDateRange.Text = $newDate
            }
        }
    }
}
